# Apply the changes described by the commit:
#  - Required Prompts cell (D2) gets its prompt names wrapped in square brackets
#  - The active selection on the sheet moves to J2
#  - Column D (Required Prompts) widens slightly to fit the new, longer text
#  - The workbook window geometry is updated to the author's new window layout

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Required Prompts" value in D2: parentName,studentName -> [parentName],[studentName]
$ws.Range("D2").Value = "[parentName],[studentName]"

# Widen column D so the longer prompt text continues to fit (was 24, now ~26.16 chars)
$ws.Columns.Item(4).ColumnWidth = 25.3

# Move the active cell/selection to J2
$null = $ws.Range("J2").Select()

# Reposition/resize the workbook window to match the author's saved layout
$win = $excel.ActiveWindow
$win.Left = 0
$win.Top = 460
$win.Width = 33600
$win.Height = 20540
